# Scheduled market-data refresh: updates the cached currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ
# columns (H:L) and the derived LeveProfitNQ / LeveProfitHQ columns (M:N)
# for the leves whose market prices changed, one crafting-job sheet at a
# time. Some rows' profit cells flip between "positive -> remove" and
# "absent -> add" because a NQ/HQ side becomes (non-)viable.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 3117.1667
$ws.Range("I98").Value = 3117.1667
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 3117.1667
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -1619.1667

# Row 122
$ws.Range("H122").Value = 3117.1667
$ws.Range("I122").Value = 3117.1667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9351.500100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -6901.500100000001

# Row 131
$ws.Range("H131").Value = 1178
$ws.Range("I131").Value = 1022.5
$ws.Range("K131").Value = 3067.5
$ws.Range("M131").Value = 1972.5

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 2260
$ws.Range("I63").Value = 1887.8572
$ws.Range("J63").Value = 3128.3333
$ws.Range("K63").Value = 1887.8572
$ws.Range("L63").Value = 3128.3333
$ws.Range("M63").Value = -1201.8572
$ws.Range("N63").Value = -4500.3333

# Row 66
$ws.Range("H66").Value = 2260
$ws.Range("I66").Value = 1887.8572
$ws.Range("J66").Value = 3128.3333
$ws.Range("K66").Value = 9439.286
$ws.Range("L66").Value = 15641.6665
$ws.Range("M66").Value = -6007.286
$ws.Range("N66").Value = -22505.6665

# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0

# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0

# Row 95
$ws.Range("H95").Value = 8750
$ws.Range("J95").Value = 8750
$ws.Range("L95").Value = 8750
$ws.Range("N95").Value = -14242

# Row 110
$ws.Range("H110").Value = 49449.094
$ws.Range("I110").Value = 53969.42
$ws.Range("K110").Value = 53969.42
$ws.Range("M110").Value = -51924.42

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 665.3333
$ws.Range("I7").Value = 698.8
$ws.Range("J7").Value = 498
$ws.Range("K7").Value = 698.8
$ws.Range("L7").Value = 498
$ws.Range("M7").Value = -585.8
$ws.Range("N7").Value = -724

# Row 15
$ws.Range("H15").Value = 750
$ws.Range("I15").Value = 500
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = -273
$ws.Range("N15").Value = -1454

# Row 23
$ws.Range("H23").Value = 799.6667
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 799.6667
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 799.6667
$ws.Range("N23").Value = -1365.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 4580.222
$ws.Range("I5").Value = 1031.8572
$ws.Range("J5").Value = 16999.5
$ws.Range("K5").Value = 1031.8572
$ws.Range("L5").Value = 16999.5
$ws.Range("M5").Value = -919.8571999999999
$ws.Range("N5").Value = -17223.5

# Row 10
$ws.Range("H10").Value = 5688.8
$ws.Range("I10").Value = 1004
$ws.Range("J10").Value = 8812
$ws.Range("K10").Value = 1004
$ws.Range("L10").Value = 8812
$ws.Range("M10").Value = -865
$ws.Range("N10").Value = -9090

# Row 12
$ws.Range("H12").Value = 3040.4443
$ws.Range("I12").Value = 2670.6
$ws.Range("J12").Value = 3502.75
$ws.Range("K12").Value = 2670.6
$ws.Range("L12").Value = 3502.75
$ws.Range("M12").Value = -2500.6
$ws.Range("N12").Value = -3842.75

# Row 13
$ws.Range("H13").Value = 166.66667
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 166.66667
$ws.Range("K13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("M13").Value = 166.66667
$ws.Range("N13").Value = -444.66667

# Row 25
$ws.Range("H25").Value = 77116.94
$ws.Range("I25").Value = 17000
$ws.Range("J25").Value = 89999.14
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 89999.14
$ws.Range("M25").Value = -16826
$ws.Range("N25").Value = -90347.14

# Row 125
$ws.Range("H125").Value = 64748.75
$ws.Range("J125").Value = 64748.75
$ws.Range("L125").Value = 64748.75
$ws.Range("N125").Value = -69668.75

# Row 134
$ws.Range("H134").Value = 50001496
$ws.Range("J134").Value = 1829.6666
$ws.Range("L134").Value = 5488.9998
$ws.Range("N134").Value = -10558.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0

# Row 115
$ws.Range("H115").Value = 5299.5
$ws.Range("J115").Value = 9999
$ws.Range("L115").Value = 29997
$ws.Range("N115").Value = -32347

# Row 118
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = 0

# Row 125
$ws.Range("H125").Value = 13996
$ws.Range("J125").Value = 13996
$ws.Range("L125").Value = 41988
$ws.Range("N125").Value = -51828

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 0

# Row 9
$ws.Range("H9").Value = 3626.5
$ws.Range("I9").Value = 1502.3334
$ws.Range("J9").Value = 9999
$ws.Range("K9").Value = 1502.3334
$ws.Range("L9").Value = 9999
$ws.Range("M9").Value = -1332.3334
$ws.Range("N9").Value = -10339

# Row 43
$ws.Range("H43").Value = 2800
$ws.Range("I43").Value = 2800
$ws.Range("K43").Value = 2800
$ws.Range("M43").Value = -2649

# Row 57
$ws.Range("H57").Value = 57499.668
$ws.Range("J57").Value = 67999.60000000001
$ws.Range("L57").Value = 67999.60000000001
$ws.Range("N57").Value = -69639.60000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7478.4116
$ws.Range("I7").Value = 7239.7856
$ws.Range("K7").Value = 7239.7856
$ws.Range("M7").Value = -7127.7856

# Row 9
$ws.Range("H9").Value = 946.8889
$ws.Range("I9").Value = 1085.6666
$ws.Range("J9").Value = 669.3333
$ws.Range("K9").Value = 1085.6666
$ws.Range("L9").Value = 669.3333
$ws.Range("M9").Value = -861.6666
$ws.Range("N9").Value = -1117.3333

# Row 93
$ws.Range("H93").Value = 725.94116
$ws.Range("I93").Value = 682.6429000000001
$ws.Range("K93").Value = 682.6429000000001
$ws.Range("M93").Value = 565.3570999999999

# Row 126
$ws.Range("H126").Value = 7478.4116
$ws.Range("I126").Value = 7239.7856
$ws.Range("K126").Value = 21719.3568
$ws.Range("M126").Value = -19249.3568

$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 33333
$ws.Range("I34").Value = 33333
$ws.Range("K34").Value = 33333
$ws.Range("M34").Value = -33130

# Row 37
$ws.Range("H37").Value = 25564
$ws.Range("I37").Value = 25564
$ws.Range("K37").Value = 25564
$ws.Range("M37").Value = -25361

# Row 126
$ws.Range("H126").Value = 1569.5
$ws.Range("I126").Value = 1446.8
$ws.Range("K126").Value = 4340.4
$ws.Range("M126").Value = -1870.4

# Row 132
$ws.Range("H132").Value = 45472480
$ws.Range("I132").Value = 50009730
$ws.Range("K132").Value = 150029190
$ws.Range("M132").Value = -150026660

# Row 135
$ws.Range("H135").Value = 88806.5
$ws.Range("J135").Value = 88806.5
$ws.Range("L135").Value = 88806.5
$ws.Range("N135").Value = -98946.5
